$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1354.070613502361
$ws.Range("C2").Value = 5078.595449712195
$ws.Range("D2").Value = 6271.615963346384
